$wb = $excel.ActiveWorkbook

# --- Overview sheet: status + generate-date for both locales, plus the
#     wider zh-cn/de-de status columns that follow from the longer text. ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2017-02-28 07:58:09"
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# --- zh-cn sheet: Status + Latest Handoff Datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2017-02-28 07:57:54"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

# --- de-de sheet: Status + Latest Handoff Datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2017-02-28 07:58:09"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
